$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 221, shifting existing rows 221:233 down to 222:234
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with the new record's data
$ws.Cells.Item(221, 1).Value = 3
$ws.Cells.Item(221, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(221, 3).Value = "Coquimbo"
$ws.Cells.Item(221, 4).Value = 44516
$ws.Cells.Item(221, 5).Value = 5
$ws.Cells.Item(221, 6).Value = 100112043
$ws.Cells.Item(221, 7).Value = "Pepino ensalada"
$ws.Cells.Item(221, 8).Value = "Sin especificar"
$ws.Cells.Item(221, 9).Value = "Primera"
$ws.Cells.Item(221, 10).Value = 125
$ws.Cells.Item(221, 11).Value = 7000
$ws.Cells.Item(221, 12).Value = 7500
$ws.Cells.Item(221, 13).Value = 7240
$ws.Cells.Item(221, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(221, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(221, 16).Value = 103
$ws.Cells.Item(221, 17).Value = 70
$ws.Cells.Item(221, 18).Value = "Hortaliza"
